$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PatientDetails")

$ws.Range("B2").Value = "Miya"
$ws.Range("C2").Value = "hanzo"
$ws.Range("D2").Value = "miyaaa@gmail.com"
$ws.Range("H2").Value = "Male"

[void]$ws.Range("D3").Select()
